$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite rows 6 and 9-18 to their final Product Backlog state ---
# (New backlog items "Implement Tutorial Elements", "Iphrit Sprite",
#  "Start Sequence for first level" and "Villagers running" were inserted;
#  existing rows 9-14 shifted down to 13-18.)

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = "Implement Tutorial Elements"
$ws.Range("D6").Value = "!!"
$ws.Range("E6").Value = "As a player I want to understand how to play the game. Furthermore I want to leran all about the mechanics of the game and how to use them."
$ws.Range("F6").ClearContents()

# Row 9
$ws.Range("A9").Value = 14
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Iphrit Sprite"
$ws.Range("D9").Value = "!"
$ws.Range("E9").Value = "As a player I want to see a nice looking Sprite of the second boss."
$ws.Range("F9").ClearContents()

# Row 10
$ws.Range("A10").Value = 15
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Start Sequence for first level"
$ws.Range("D10").Value = "!!"
$ws.Range("E10").Value = "As a player I want to see a sequence why the player gets out of his hut, to investigate whats going on."
$ws.Range("F10").ClearContents()

# Row 11
$ws.Range("A11").Value = 16
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "Villagers running "
$ws.Range("D11").Value = "!"
$ws.Range("E11").Value = "As a player I want to see the villagers running by the hut of the Adventurer. Furthermore I want to hear an explanation why they are running away."
$ws.Range("F11").ClearContents()

# Row 12
$ws.Range("A12").Value = 17
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()

# Row 13
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Design new Enemies"
$ws.Range("D13").Value = "!"
$ws.Range("E13").Value = "As a Player I want to fight against 100  enemies. That means 10 new enemies per level. I also want to see some familiar enemies which I already met in previous games."
$ws.Range("F13").Value = 5

# Row 14
$ws.Range("A14").Value = 9
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Convert Level form Paper to Unity"
$ws.Range("D14").Value = "!"
$ws.Range("E14").Value = "As a player I want to have a good structured level, where I can move around freely and fight enemies."
$ws.Range("F14").Value = 5

# Row 15
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Design Dialoges"
$ws.Range("D15").Value = "!"
$ws.Range("E15").Value = "As a player I want to have interesting and hilarious dialoges between the characters.  "
$ws.Range("F15").Value = 2

# Row 16
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Implement Questlog"
$ws.Range("D16").Value = "!"
$ws.Range("E16").Value = "As a player I want to have a widget to access all my quests in order to have a good overview of all my open quests. Furthermore I want to get informed if a new quest comes up (Display Box for ""You got a new Quest"")"
$ws.Range("F16").Value = 5

# Row 17
$ws.Range("A17").Value = 12
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "Old Lady Quest"
$ws.Range("D17").Value = "!!"
$ws.Range("E17").Value = "As a player I want to get  the quest to help the old lady bring her ""groceries"" to her hut."
$ws.Range("F17").Value = 5

# Row 18
$ws.Range("A18").Value = 13
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "Hut of the old Lady"
$ws.Range("D18").Value = "!!"
$ws.Range("E18").Value = "As a player I want to see the hut of the old Lady. Furthermore I want to move around in it."
$ws.Range("F18").Value = 2

# --- Row heights: rows whose wrapped "Vertics" text now spans two lines ---
$ws.Rows("6:6").RowHeight = 30
$ws.Rows("9:9").RowHeight = 15
$ws.Rows("11:11").RowHeight = 30
$ws.Rows("12:12").RowHeight = 15
$ws.Rows("13:13").RowHeight = 30
$ws.Rows("16:16").RowHeight = 30

# --- Selection moved to C12 (the new empty backlog-item placeholder row) ---
$ws.Range("C12").Select()

